# Macroferia Regional de Talca - Mandarina: insert a new weekly price record.
# A new row is inserted at row 167 (pushing the existing rows 167-244 down to
# 168-245) and populated with the latest report for Murcott / Primera.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row above the current row 167 - this shifts every
# existing row (167..244) down by one (to 168..245) and keeps all of their
# data/formatting intact.
$ws.Rows.Item(167).Insert()

# Populate the newly inserted row 167 with the new record.
$ws.Cells.Item(167, 1).Value = 5
$ws.Cells.Item(167, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(167, 3).Value = "Maule"
$ws.Cells.Item(167, 4).Value = 44572
$ws.Cells.Item(167, 5).Value = 7
$ws.Cells.Item(167, 6).Value = "Fruta"
$ws.Cells.Item(167, 7).Value = 100102
$ws.Cells.Item(167, 8).Value = "Cítricos"
$ws.Cells.Item(167, 9).Value = 100102004
$ws.Cells.Item(167, 10).Value = "Mandarina"
$ws.Cells.Item(167, 11).Value = "Murcott"
$ws.Cells.Item(167, 12).Value = "Primera"
$ws.Cells.Item(167, 13).Value = 500
$ws.Cells.Item(167, 14).Value = 7000
$ws.Cells.Item(167, 15).Value = 7000
$ws.Cells.Item(167, 16).Value = 7000
$ws.Cells.Item(167, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(167, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(167, 19).Value = 389
$ws.Cells.Item(167, 20).Value = 18
